$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (all values are stored as text
# in the source workbook, matching the original inline-string formatting)
$updates = [ordered]@{
    "D2" = "294.73"
    "E2" = "2.01%"
    "D3" = "31.04"
    "E3" = "0.01%"
    "D4" = "4.928"
    "E4" = "0.35%"
    "D5" = "0.07370"
    "E5" = "3.43%"
    "D6" = "2.312"
    "E6" = "29.01%"
    "D7" = "7.715"
    "E7" = "0.86%"
    "D8" = "3.757"
    "E8" = "-0.59%"
    "D9" = "0.9129"
    "E9" = "2.09%"
    "E10" = "2.65%"
    "D11" = "0.08259"
    "E11" = "9.57%"
    "D12" = "0.08254"
    "E12" = "2.16%"
    "D13" = "0.03119"
    "E13" = "4.34%"
    "D14" = "0.1007"
    "E14" = "0.83%"
    "D15" = "0.001508"
    "E15" = "0.76%"
    "D16" = "0.005741"
    "E16" = "-0.52%"
    "E17" = "0.53%"
    "D18" = "2.079"
    "E18" = "-1.32%"
    "E19" = "1.63%"
    "E20" = "0.48%"
    "D21" = "3.964"
    "E21" = "-7.20%"
    "D22" = "0.2099"
    "E22" = "4.25%"
    "D23" = "0.04550"
    "E23" = "1.43%"
    "D24" = "0.001210"
    "E24" = "-0.19%"
    "D25" = "0.004341"
    "E25" = "-6.77%"
    "D26" = "0.0001300"
    "E26" = "3.74%"
    "D27" = "0.0003393"
    "D39" = "0.01608"
    "E39" = "-1.61%"
    "D40" = "0.04444"
    "D41" = "0.007335"
    "E41" = "-0.79%"
    "D42" = "0.008818"
    "D43" = "0.1325"
    "E43" = "1.56%"
    "D44" = "0.002069"
    "E44" = "3.18%"
    "D45" = "0.009104"
    "E45" = "-11.88%"
    "D46" = "0.00006041"
    "E46" = "2.96%"
    "E47" = "-0.10%"
    "E48" = "1.05%"
    "D50" = "0.00002099"
    "E50" = "-0.10%"
    "D51" = "0.0001999"
    "E51" = "-0.10%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text format so Excel keeps the literal string instead of
    # auto-converting numeric-looking / percent-looking text to a number
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}

Write-Output "Updated $($updates.Count) cells"
